# mis à jour de la liste de taches
#
# Rows 47, 48, 62, 63: task was "En cours" (col C) and is now "Terminé" (col D).
# Move the "X" mark from column C to column D.
#
# Rows 50, 51: tasks were still "A faire" (only col B ticked, no responsible /
# dates). They are now "En cours" (col C ticked) with a responsible
# ("Florentin") and a start date, matching the same look (fill/format) as the
# other "En cours" rows (47, 48, 62, 63).
#
# Also update the current selection to D67 (was F63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 47 / 48 / 62 / 63: move the "X" from column C to column D ---
foreach ($r in 47, 48, 62, 63) {
    $ws.Range("C$r").ClearContents()
    $ws.Range("D$r").Value = "X"
}

# --- Rows 50 / 51: re-format to match the "En cours" rows, then set values ---
foreach ($r in 50, 51) {
    $ws.Range("B47:H47").Copy()
    $ws.Range("B$r`:H$r").PasteSpecial(-4122)

    $ws.Range("B$r").ClearContents()
    $ws.Range("C$r").Value = "X"
    $ws.Range("D$r").ClearContents()
    $ws.Range("E$r").Value = "Florentin"
    $ws.Range("F$r").Value = 42689
    $ws.Range("G$r").ClearContents()
    $ws.Range("H$r").ClearContents()
}

$excel.CutCopyMode = 0

# --- Update the selection shown when the sheet is reopened ---
$ws.Range("D67").Select() | Out-Null
